# Add a new "Sede Reason" column (O) to the inquiring requests report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (N1, "Domain Reason")
# onto the new header cell (O1) so it keeps the same shaded header style,
# then set its text.
$ws.Cells.Item(1, 14).Copy()
$ws.Cells.Item(1, 15).PasteSpecial(-4122)
$ws.Cells.Item(1, 15).Value = "Sede Reason"

# Give the new column a sensible width matching its sibling "reason" columns.
$ws.Columns.Item(15).ColumnWidth = 16.140625

# Re-apply the autofilter over the extended range A1:O1.
$ws.AutoFilterMode = $false
$ws.Range("A1:O1").AutoFilter(1)

# Keep the hidden _FilterDatabase defined name in sync with the new filter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$O`$1"
    }
}

# Restore the cursor/selection position recorded in the saved workbook.
$ws.Range("I14").Select()
